$d = $word.ActiveDocument

# --- Step 1: fix up the last data row of the table (row 17), which in the
# edited document drops its old 3rd/4th cell values and gains two brand new
# trailing cells. Since the cell count stays the same (5), we just overwrite
# each of the 5 cells' text directly, cell by cell, to avoid any interaction
# with the global Find/Replace pass below.
$t = $d.Tables.Item(1)
$row = $t.Rows.Item(17)
$row.Cells.Item(1).Range.Text = "45÷7=6, 3"
$row.Cells.Item(2).Range.Text = "80÷8=10, 0"
$row.Cells.Item(3).Range.Text = "91÷8=11, 3"
$row.Cells.Item(4).Range.Text = "77÷5=15, 2"
$row.Cells.Item(5).Range.Text = "16÷6=2, 4"

# --- Step 2: straightforward whole-number-problem text swaps elsewhere in
# the table. Each "old" value is unique in the document, so a global
# Find/Replace (wdReplaceAll) is safe.
$replacements = @(
    @("45÷9=5, 0", "50÷6=8, 2"),
    @("81÷6=13, 3", "40÷7=5, 5"),
    @("34÷6=5, 4", "39÷2=19, 1"),
    @("64÷7=9, 1", "31÷4=7, 3"),
    @("43÷9=4, 7", "27÷9=3, 0"),
    @("54÷4=13, 2", "20÷2=10, 0"),
    @("34÷3=11, 1", "30÷7=4, 2"),
    @("46÷3=15, 1", "75÷6=12, 3"),
    @("51÷2=25, 1", "15÷2=7, 1"),
    @("84÷9=9, 3", "81÷9=9, 0"),
    @("63÷4=15, 3", "70÷6=11, 4"),
    @("21÷7=3, 0", "58÷9=6, 4"),
    @("60÷8=7, 4", "71÷6=11, 5"),
    @("35÷2=17, 1", "96÷6=16, 0"),
    @("44÷4=11, 0", "70÷6=11, 4"),
    @("42÷9=4, 6", "73÷5=14, 3"),
    @("34÷2=17, 0", "58÷6=9, 4"),
    @("90÷9=10, 0", "89÷9=9, 8"),
    @("94÷7=13, 3", "73÷6=12, 1"),
    @("52÷4=13, 0", "92÷8=11, 4")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Write-Host "Done applying edits"
